$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-label the angle groups (-30/-15/0/15/30 -> -60/-30/0/30/60) ---
# Update the three cells whose text already carries a "quote prefix" style
# (A2, A6, A10) using a leading apostrophe so Excel keeps treating the
# label as literal text (preserves style s="1"). A14/A18 never needed the
# quote-prefix style, so they are set normally (keeps style s="4").
$ws.Range("A6").Value = "'-30°"
$ws.Range("A14").Value = "30°"
$ws.Range("A2").Value = "'-60°"
$ws.Range("A18").Value = "60°"
$ws.Range("A10").Value = "'0°"

# --- Header resistor label: 300kΩ -> 20MΩ (two-run rich text, 2nd run Arial) ---
$ws.Range("A1").Value = "20MΩ"
$ws.Range("A1").Characters(4, 1).Font.Name = "Arial"

# --- Fill in the measured readings (previously blank) ---
# Group -60° (row 2-5)
$ws.Range("C2").Value = 976.44
$ws.Range("D2").Value = 980.76
$ws.Range("E2").Value = 992.2
$ws.Range("F2").Value = 995.1

$ws.Range("C3").Value = 946.5
$ws.Range("D3").Value = 961.28
$ws.Range("E3").Value = 976.32
$ws.Range("F3").Value = 983.66

$ws.Range("C4").Value = 994.4
$ws.Range("D4").Value = 1000.54
$ws.Range("E4").Value = 1006.28
$ws.Range("F4").Value = 1008.06

$ws.Range("C5").Value = 1020.36
$ws.Range("D5").Value = 1019.72
$ws.Range("E5").Value = 1018.48
$ws.Range("F5").Value = 1020.28

# Group -30° (row 6-9)
$ws.Range("C6").Value = 950.08
$ws.Range("D6").Value = 975.46
$ws.Range("E6").Value = 987.28
$ws.Range("F6").Value = 997

$ws.Range("C7").Value = 846.12
$ws.Range("D7").Value = 876.18
$ws.Range("E7").Value = 882.32
$ws.Range("F7").Value = 919.8

$ws.Range("C8").Value = 958.62
$ws.Range("D8").Value = 976.68
$ws.Range("E8").Value = 968.92
$ws.Range("F8").Value = 983.66

$ws.Range("C9").Value = 1019.2
$ws.Range("D9").Value = 1020.02
$ws.Range("E9").Value = 1019.78
$ws.Range("F9").Value = 1018.36

# Group 0° (row 10-13)
$ws.Range("C10").Value = 966.74
$ws.Range("D10").Value = 978.36
$ws.Range("E10").Value = 988.44
$ws.Range("F10").Value = 994.98

$ws.Range("C11").Value = 141.44
$ws.Range("D11").Value = 147.62
$ws.Range("E11").Value = 271.3
$ws.Range("F11").Value = 468.16

$ws.Range("C12").Value = 744.08
$ws.Range("D12").Value = 744.26
$ws.Range("E12").Value = 748.54
$ws.Range("F12").Value = 816.14

$ws.Range("C13").Value = 1020.42
$ws.Range("D13").Value = 1020.06
$ws.Range("E13").Value = 1020.54
$ws.Range("F13").Value = 1017.92

# Group 30° (row 14-17)
$ws.Range("C14").Value = 968.12
$ws.Range("D14").Value = 980.2
$ws.Range("E14").Value = 989.5
$ws.Range("F14").Value = 995.2

$ws.Range("C15").Value = 118.24
$ws.Range("D15").Value = 127.68
$ws.Range("E15").Value = 140.18
$ws.Range("F15").Value = 151.28

$ws.Range("C16").Value = 744.82
$ws.Range("D16").Value = 745.4
$ws.Range("E16").Value = 744.48
$ws.Range("F16").Value = 743.96

$ws.Range("C17").Value = 1019.44
$ws.Range("D17").Value = 1019
$ws.Range("E17").Value = 1020.02
$ws.Range("F17").Value = 1019.96

# Group 60° (row 18-21)
$ws.Range("C18").Value = 976.24
$ws.Range("D18").Value = 976.28
$ws.Range("E18").Value = 987.82
$ws.Range("F18").Value = 995.16

$ws.Range("C19").Value = 112.46
$ws.Range("D19").Value = 123.68
$ws.Range("E19").Value = 138.38
$ws.Range("F19").Value = 150.5

$ws.Range("C20").Value = 744.7
$ws.Range("D20").Value = 744.62
$ws.Range("E20").Value = 744.52
$ws.Range("F20").Value = 744.24

$ws.Range("C21").Value = 1020.48
$ws.Range("D21").Value = 1020.54
$ws.Range("E21").Value = 1020.32
$ws.Range("F21").Value = 1020.32

# --- Notes row ---
$ws.Range("A23").Value = "PT1 & PT2 seem okay, PT3 & PT4 less so"

# --- Match the saved selection state (A24) ---
$ws.Range("A24").Select()
